# AutoSPInstaller Issue Tracker: add 'Tested with 3.87' solution note to row 6,
# and append a new row 8 documenting the SP2013 Access Services / mixed-mode
# authentication issue and its work-around.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: existing "Managed Metadata Service" issue gets a Solution note ---
# E6 (Work around) switches from wrap-only style to vertical-top + wrap style
$ws.Range("E6").VerticalAlignment = -4160

# F6 (Solution) is new text
$ws.Range("F6").Value2 = "Tested with AutoSPInstaller Version 3.87, and issue appears to no longer occur."
$ws.Range("F6").VerticalAlignment = -4160
$ws.Range("F6").WrapText = $true

# --- New Row 8: SP2013 Access Services / mixed mode authentication issue ---
$ws.Range("A8").Value2 = 2013
$ws.Range("A8").VerticalAlignment = -4160

$ws.Range("B8").Value2 = 3.87
$ws.Range("B8").VerticalAlignment = -4160

$ws.Range("F8").Value2 = "Ensure that your SQL Server instance has mixed mode authentication enabled for both Windows and SQL logins. Properties --> Security --> SQL Server and Windows Authentication mode."
$ws.Range("F8").VerticalAlignment = -4160
$ws.Range("F8").WrapText = $true

$ws.Range("D8").Value2 = "When you select 'True' to provision Access Services 2013 under your Enterprise Service Applications; and then encounter an exception message like the one below when you run the 'AutoSPInstallerLaunch' BAT file. 'New-SPAccessServicesApplication : A connection could be established to the Application Database Server but mixed mode authentication isn't enabled.'"
$ws.Range("D8").VerticalAlignment = -4160
$ws.Range("D8").WrapText = $true

$ws.Rows.Item(8).RowHeight = 210

# Selection ends up on the newly added row
$ws.Range("A8").Select()
